$wb = $excel.ActiveWorkbook

# --- Sheet "37alt": append day3 raw OD rows (50-73) ---
$ws1 = $wb.Worksheets.Item("37alt")

# Clone the formatting of existing data rows onto the new block before
# writing values, so the appended rows keep the same font/alignment/fill
# as the rest of the table (style group "A3:D3" for A-col style 3 wells,
# style group "A14:D14" for A-col style 2 wells).
$ws1.Range("A3:D3").Copy()
$ws1.Range("A50:D61").PasteSpecial(-4122)
$ws1.Range("A14:D14").Copy()
$ws1.Range("A62:D73").PasteSpecial(-4122)

$ws1.Range("A50").Value = "B2"
$ws1.Range("B50").Value = "culture"
$ws1.Range("C50").Value = 0.22579999268054962
$ws1.Range("D50").Value = "day3"
$ws1.Range("A51").Value = "B4"
$ws1.Range("B51").Value = "blank"
$ws1.Range("C51").Value = 0.05000000074505806
$ws1.Range("D51").Value = "day3"
$ws1.Range("A52").Value = "B6"
$ws1.Range("B52").Value = "culture"
$ws1.Range("C52").Value = 0.24789999425411224
$ws1.Range("D52").Value = "day3"
$ws1.Range("A53").Value = "B8"
$ws1.Range("B53").Value = "culture"
$ws1.Range("C53").Value = 0.15620000660419464
$ws1.Range("D53").Value = "day3"
$ws1.Range("A54").Value = "C5"
$ws1.Range("B54").Value = "culture"
$ws1.Range("C54").Value = 0.2110999971628189
$ws1.Range("D54").Value = "day3"
$ws1.Range("A55").Value = "C7"
$ws1.Range("B55").Value = "culture"
$ws1.Range("C55").Value = 0.20970000326633453
$ws1.Range("D55").Value = "day3"
$ws1.Range("A56").Value = "C9"
$ws1.Range("B56").Value = "culture"
$ws1.Range("C56").Value = 0.2134000062942505
$ws1.Range("D56").Value = "day3"
$ws1.Range("A57").Value = "C11"
$ws1.Range("B57").Value = "blank"
$ws1.Range("C57").Value = 0.0478999987244606
$ws1.Range("D57").Value = "day3"
$ws1.Range("A58").Value = "D2"
$ws1.Range("B58").Value = "culture"
$ws1.Range("C58").Value = 0.23469999432563782
$ws1.Range("D58").Value = "day3"
$ws1.Range("A59").Value = "D4"
$ws1.Range("B59").Value = "culture"
$ws1.Range("C59").Value = 0.3447999954223633
$ws1.Range("D59").Value = "day3"
$ws1.Range("A60").Value = "D6"
$ws1.Range("B60").Value = "blank"
$ws1.Range("C60").Value = 0.04899999871850014
$ws1.Range("D60").Value = "day3"
$ws1.Range("A61").Value = "D8"
$ws1.Range("B61").Value = "culture"
$ws1.Range("C61").Value = 0.3109000027179718
$ws1.Range("D61").Value = "day3"
$ws1.Range("A62").Value = "E5"
$ws1.Range("B62").Value = "blank"
$ws1.Range("C62").Value = 0.05400000140070915
$ws1.Range("D62").Value = "day3"
$ws1.Range("A63").Value = "E7"
$ws1.Range("B63").Value = "culture"
$ws1.Range("C63").Value = 0.17730000615119934
$ws1.Range("D63").Value = "day3"
$ws1.Range("A64").Value = "E9"
$ws1.Range("B64").Value = "culture"
$ws1.Range("C64").Value = 0.20909999310970306
$ws1.Range("D64").Value = "day3"
$ws1.Range("A65").Value = "E11"
$ws1.Range("B65").Value = "culture"
$ws1.Range("C65").Value = 0.163100004196167
$ws1.Range("D65").Value = "day3"
$ws1.Range("A66").Value = "F2"
$ws1.Range("B66").Value = "culture"
$ws1.Range("C66").Value = 0.2069000005722046
$ws1.Range("D66").Value = "day3"
$ws1.Range("A67").Value = "F4"
$ws1.Range("B67").Value = "culture"
$ws1.Range("C67").Value = 0.22360000014305115
$ws1.Range("D67").Value = "day3"
$ws1.Range("A68").Value = "F6"
$ws1.Range("B68").Value = "culture"
$ws1.Range("C68").Value = 0.25110000371932983
$ws1.Range("D68").Value = "day3"
$ws1.Range("A69").Value = "F8"
$ws1.Range("B69").Value = "blank"
$ws1.Range("C69").Value = 0.05090000107884407
$ws1.Range("D69").Value = "day3"
$ws1.Range("A70").Value = "G5"
$ws1.Range("B70").Value = "culture"
$ws1.Range("C70").Value = 0.28380000591278076
$ws1.Range("D70").Value = "day3"
$ws1.Range("A71").Value = "G7"
$ws1.Range("B71").Value = "blank"
$ws1.Range("C71").Value = 0.05050000175833702
$ws1.Range("D71").Value = "day3"
$ws1.Range("A72").Value = "G9"
$ws1.Range("B72").Value = "culture"
$ws1.Range("C72").Value = 0.2547999918460846
$ws1.Range("D72").Value = "day3"
$ws1.Range("A73").Value = "G11"
$ws1.Range("B73").Value = "culture"
$ws1.Range("C73").Value = 0.22859999537467957
$ws1.Range("D73").Value = "day3"

# --- Sheet "42alt": append day3 raw OD rows (50-73) ---
$ws2 = $wb.Worksheets.Item("42alt")

$ws2.Range("A3:D3").Copy()
$ws2.Range("A50:D61").PasteSpecial(-4122)
$ws2.Range("A14:D14").Copy()
$ws2.Range("A62:D73").PasteSpecial(-4122)

$ws2.Range("A50").Value = "B2"
$ws2.Range("B50").Value = "culture"
$ws2.Range("C50").Value = 0.05220000073313713
$ws2.Range("D50").Value = "day3"
$ws2.Range("A51").Value = "B4"
$ws2.Range("B51").Value = "blank"
$ws2.Range("C51").Value = 0.051500000059604645
$ws2.Range("D51").Value = "day3"
$ws2.Range("A52").Value = "B6"
$ws2.Range("B52").Value = "culture"
$ws2.Range("C52").Value = 0.050599999725818634
$ws2.Range("D52").Value = "day3"
$ws2.Range("A53").Value = "B8"
$ws2.Range("B53").Value = "culture"
$ws2.Range("C53").Value = 0.050200000405311584
$ws2.Range("D53").Value = "day3"
$ws2.Range("A54").Value = "C5"
$ws2.Range("B54").Value = "culture"
$ws2.Range("C54").Value = 0.04969999939203262
$ws2.Range("D54").Value = "day3"
$ws2.Range("A55").Value = "C7"
$ws2.Range("B55").Value = "culture"
$ws2.Range("C55").Value = 0.0478999987244606
$ws2.Range("D55").Value = "day3"
$ws2.Range("A56").Value = "C9"
$ws2.Range("B56").Value = "culture"
$ws2.Range("C56").Value = 0.050200000405311584
$ws2.Range("D56").Value = "day3"
$ws2.Range("A57").Value = "C11"
$ws2.Range("B57").Value = "blank"
$ws2.Range("C57").Value = 0.051600001752376556
$ws2.Range("D57").Value = "day3"
$ws2.Range("A58").Value = "D2"
$ws2.Range("B58").Value = "culture"
$ws2.Range("C58").Value = 0.0502999983727932
$ws2.Range("D58").Value = "day3"
$ws2.Range("A59").Value = "D4"
$ws2.Range("B59").Value = "culture"
$ws2.Range("C59").Value = 0.048900000751018524
$ws2.Range("D59").Value = "day3"
$ws2.Range("A60").Value = "D6"
$ws2.Range("B60").Value = "blank"
$ws2.Range("C60").Value = 0.04820000007748604
$ws2.Range("D60").Value = "day3"
$ws2.Range("A61").Value = "D8"
$ws2.Range("B61").Value = "culture"
$ws2.Range("C61").Value = 0.051100000739097595
$ws2.Range("D61").Value = "day3"
$ws2.Range("A62").Value = "E5"
$ws2.Range("B62").Value = "blank"
$ws2.Range("C62").Value = 0.04919999837875366
$ws2.Range("D62").Value = "day3"
$ws2.Range("A63").Value = "E7"
$ws2.Range("B63").Value = "culture"
$ws2.Range("C63").Value = 0.049300000071525574
$ws2.Range("D63").Value = "day3"
$ws2.Range("A64").Value = "E9"
$ws2.Range("B64").Value = "culture"
$ws2.Range("C64").Value = 0.04910000041127205
$ws2.Range("D64").Value = "day3"
$ws2.Range("A65").Value = "E11"
$ws2.Range("B65").Value = "culture"
$ws2.Range("C65").Value = 0.04989999905228615
$ws2.Range("D65").Value = "day3"
$ws2.Range("A66").Value = "F2"
$ws2.Range("B66").Value = "culture"
$ws2.Range("C66").Value = 0.05119999870657921
$ws2.Range("D66").Value = "day3"
$ws2.Range("A67").Value = "F4"
$ws2.Range("B67").Value = "culture"
$ws2.Range("C67").Value = 0.050200000405311584
$ws2.Range("D67").Value = "day3"
$ws2.Range("A68").Value = "F6"
$ws2.Range("B68").Value = "culture"
$ws2.Range("C68").Value = 0.048700001090765
$ws2.Range("D68").Value = "day3"
$ws2.Range("A69").Value = "F8"
$ws2.Range("B69").Value = "blank"
$ws2.Range("C69").Value = 0.0494999997317791
$ws2.Range("D69").Value = "day3"
$ws2.Range("A70").Value = "G5"
$ws2.Range("B70").Value = "culture"
$ws2.Range("C70").Value = 0.04910000041127205
$ws2.Range("D70").Value = "day3"
$ws2.Range("A71").Value = "G7"
$ws2.Range("B71").Value = "blank"
$ws2.Range("C71").Value = 0.04989999905228615
$ws2.Range("D71").Value = "day3"
$ws2.Range("A72").Value = "G9"
$ws2.Range("B72").Value = "culture"
$ws2.Range("C72").Value = 0.050599999725818634
$ws2.Range("D72").Value = "day3"
$ws2.Range("A73").Value = "G11"
$ws2.Range("B73").Value = "culture"
$ws2.Range("C73").Value = 0.05350000038743019
$ws2.Range("D73").Value = "day3"

